$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 5-9 columns B, D, F (these rows will only keep column A)
$ws.Range("B5:B9").ClearContents()
$ws.Range("D5:D9").ClearContents()
$ws.Range("F5:F9").ClearContents()

# Update row 10
$ws.Range("B10").Value = "Didier"
$ws.Range("D10").Value = "Mathieu"
$ws.Range("F10").Value = "07:52"

# Add row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Hugo D"
$ws.Range("D11").Value = "Didier"
$ws.Range("F11").Value = "07:48"

# Add row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Anne-Lise"
$ws.Range("D12").Value = "Didier"
$ws.Range("F12").Value = "07:47"
